$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that flip from "dimension" to "measure": A, D, E, H, I, K
$cols = @("A", "D", "E", "H", "I", "K")

foreach ($col in $cols) {
    # Row 2: iaest-dimension:<x> -> iaest-measure:<x>
    $cell2 = $ws.Range($col + "2")
    $oldVal = $cell2.Value()
    $cell2.Value = $oldVal -replace "^iaest-dimension:", "iaest-measure:"

    # Row 3: dim -> medida
    $ws.Range($col + "3").Value = "medida"

    # Row 4: skos:Concept -> xsd:int
    $ws.Range($col + "4").Value = "xsd:int"

    # Row 5: remove the mapping-*.xlsx reference cell entirely
    $ws.Range($col + "5").Clear()
}
